# Update Leve profit calculations across the Sheets workbook
# per scheduled market-data refresh (Halicarnassus_Profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 904.2857
$ws.Cells.Item(19, 10).Value = 696.5
$ws.Cells.Item(19, 12).Value = 696.5
$ws.Cells.Item(19, 14).Value = -1046.5
$ws.Cells.Item(70, 8).Value = 3633.9412
$ws.Cells.Item(70, 9).Value = 1730.6666
$ws.Cells.Item(70, 10).Value = 4041.7856
$ws.Cells.Item(70, 11).Value = 5191.9998
$ws.Cells.Item(70, 12).Value = 12125.3568
$ws.Cells.Item(70, 13).Value = -4921.9998
$ws.Cells.Item(70, 14).Value = -12665.3568
$ws.Cells.Item(73, 8).Value = 3633.9412
$ws.Cells.Item(73, 9).Value = 1730.6666
$ws.Cells.Item(73, 10).Value = 4041.7856
$ws.Cells.Item(73, 11).Value = 5191.9998
$ws.Cells.Item(73, 12).Value = 12125.3568
$ws.Cells.Item(73, 13).Value = -4255.9998
$ws.Cells.Item(73, 14).Value = -13997.3568
$ws.Cells.Item(80, 8).Value = 615
$ws.Cells.Item(80, 9).Value = 550
$ws.Cells.Item(80, 10).Value = 631.25
$ws.Cells.Item(80, 11).Value = 1650
$ws.Cells.Item(80, 12).Value = 1893.75
$ws.Cells.Item(80, 13).Value = -652
$ws.Cells.Item(80, 14).Value = -3889.75
$ws.Cells.Item(83, 8).Value = 615
$ws.Cells.Item(83, 9).Value = 550
$ws.Cells.Item(83, 10).Value = 631.25
$ws.Cells.Item(83, 11).Value = 4950
$ws.Cells.Item(83, 12).Value = 5681.25
$ws.Cells.Item(83, 13).Value = 42
$ws.Cells.Item(83, 14).Value = -15665.25
$ws.Cells.Item(135, 8).Value = 1103.4286
$ws.Cells.Item(135, 9).Value = 1103.4286
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 9930.857399999999
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).ClearContents()
$ws.Cells.Item(135, 14).Value = -7395.857399999999
$ws.Cells.Item(137, 8).Value = 1396.4445
$ws.Cells.Item(137, 10).Value = 2796
$ws.Cells.Item(137, 12).Value = 8388
$ws.Cells.Item(137, 14).Value = -13488
$ws.Cells.Item(138, 8).Value = 3040.1667
$ws.Cells.Item(138, 9).Value = 847.5454999999999
$ws.Cells.Item(138, 10).Value = 4895.4614
$ws.Cells.Item(138, 11).Value = 2542.6365
$ws.Cells.Item(138, 12).Value = 14686.3842
$ws.Cells.Item(138, 13).Value = 2597.3635
$ws.Cells.Item(138, 14).Value = -24966.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 3982.6365
$ws.Cells.Item(102, 9).Value = 1618.3334
$ws.Cells.Item(102, 10).Value = 6819.8
$ws.Cells.Item(102, 11).Value = 1618.3334
$ws.Cells.Item(102, 12).Value = 6819.8
$ws.Cells.Item(102, 13).Value = 3.666600000000017
$ws.Cells.Item(102, 14).Value = -10063.8
$ws.Cells.Item(104, 8).Value = 504
$ws.Cells.Item(104, 9).Value = 504
$ws.Cells.Item(104, 11).Value = 504
$ws.Cells.Item(104, 13).Value = 2990
$ws.Cells.Item(132, 8).Value = 3370.8
$ws.Cells.Item(132, 9).Value = 3338.147
$ws.Cells.Item(132, 11).Value = 10014.441
$ws.Cells.Item(132, 13).Value = -7484.440999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).ClearContents()
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = 0
$ws.Cells.Item(107, 8).Value = 2540.125
$ws.Cells.Item(107, 10).Value = 8503.25
$ws.Cells.Item(107, 12).Value = 8503.25
$ws.Cells.Item(107, 14).Value = -12343.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 4089.1155
$ws.Cells.Item(7, 9).Value = 6791.2
$ws.Cells.Item(7, 10).Value = 404.45456
$ws.Cells.Item(7, 11).Value = 6791.2
$ws.Cells.Item(7, 12).Value = 404.45456
$ws.Cells.Item(7, 13).Value = -6678.2
$ws.Cells.Item(7, 14).Value = -630.45456
$ws.Cells.Item(22, 8).Value = 2499.875
$ws.Cells.Item(22, 9).Value = 2416.1667
$ws.Cells.Item(22, 10).Value = 2751
$ws.Cells.Item(22, 11).Value = 2416.1667
$ws.Cells.Item(22, 12).Value = 2751
$ws.Cells.Item(22, 13).Value = -2066.1667
$ws.Cells.Item(22, 14).Value = -3451
$ws.Cells.Item(31, 8).Value = 4024
$ws.Cells.Item(31, 9).Value = 1988.7407
$ws.Cells.Item(31, 10).Value = 8251.076999999999
$ws.Cells.Item(31, 11).Value = 1988.7407
$ws.Cells.Item(31, 12).Value = 8251.076999999999
$ws.Cells.Item(31, 13).Value = -1693.7407
$ws.Cells.Item(31, 14).Value = -8841.076999999999
$ws.Cells.Item(34, 8).Value = 4024
$ws.Cells.Item(34, 9).Value = 1988.7407
$ws.Cells.Item(34, 10).Value = 8251.076999999999
$ws.Cells.Item(34, 11).Value = 1988.7407
$ws.Cells.Item(34, 12).Value = 8251.076999999999
$ws.Cells.Item(34, 13).Value = -1786.7407
$ws.Cells.Item(34, 14).Value = -8655.076999999999
$ws.Cells.Item(96, 8).Value = 8251
$ws.Cells.Item(96, 10).Value = 8251
$ws.Cells.Item(96, 12).Value = 8251
$ws.Cells.Item(96, 14).Value = -13743
$ws.Cells.Item(107, 8).Value = 350.2857
$ws.Cells.Item(107, 9).Value = 223.16667
$ws.Cells.Item(107, 10).Value = 1113
$ws.Cells.Item(107, 11).Value = 223.16667
$ws.Cells.Item(107, 12).Value = 1113
$ws.Cells.Item(107, 13).Value = 1696.83333
$ws.Cells.Item(107, 14).Value = -4953
$ws.Cells.Item(132, 8).Value = 2000
$ws.Cells.Item(132, 9).Value = 2000
$ws.Cells.Item(132, 11).Value = 6000
$ws.Cells.Item(132, 13).Value = -3470

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 950
$ws.Cells.Item(63, 9).Value = 950
$ws.Cells.Item(63, 11).Value = 2850
$ws.Cells.Item(63, 13).Value = -2101
$ws.Cells.Item(66, 8).Value = 950
$ws.Cells.Item(66, 9).Value = 950
$ws.Cells.Item(66, 11).Value = 8550
$ws.Cells.Item(66, 13).Value = -4806
$ws.Cells.Item(80, 8).Value = 4175
$ws.Cells.Item(80, 10).Value = 4569.375
$ws.Cells.Item(80, 12).Value = 13708.125
$ws.Cells.Item(80, 14).Value = -15580.125
$ws.Cells.Item(83, 8).Value = 4175
$ws.Cells.Item(83, 10).Value = 4569.375
$ws.Cells.Item(83, 12).Value = 41124.375
$ws.Cells.Item(83, 14).Value = -50484.375
$ws.Cells.Item(92, 8).Value = 539.8
$ws.Cells.Item(92, 9).Value = 600.6667
$ws.Cells.Item(92, 10).Value = 448.5
$ws.Cells.Item(92, 11).Value = 1802.0001
$ws.Cells.Item(92, 12).Value = 1345.5
$ws.Cells.Item(92, 13).Value = -554.0001
$ws.Cells.Item(92, 14).Value = -3841.5
$ws.Cells.Item(140, 8).Value = 2809
$ws.Cells.Item(140, 9).Value = 1713.5
$ws.Cells.Item(140, 11).Value = 5140.5
$ws.Cells.Item(140, 13).Value = 39.5
$ws.Cells.Item(141, 8).Value = 1000
$ws.Cells.Item(141, 9).Value = 1000
$ws.Cells.Item(141, 11).Value = 3000
$ws.Cells.Item(141, 13).Value = 2180

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4113.5713
$ws.Cells.Item(70, 9).Value = 3966
$ws.Cells.Item(70, 11).Value = 3966
$ws.Cells.Item(70, 13).Value = -3696
$ws.Cells.Item(73, 8).Value = 4113.5713
$ws.Cells.Item(73, 9).Value = 3966
$ws.Cells.Item(73, 11).Value = 3966
$ws.Cells.Item(73, 13).Value = -3030
$ws.Cells.Item(126, 8).Value = 2819
$ws.Cells.Item(126, 9).Value = 2796
$ws.Cells.Item(126, 11).Value = 8388
$ws.Cells.Item(126, 13).Value = -5918

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6506
$ws.Cells.Item(7, 9).Value = 6438.1665
$ws.Cells.Item(7, 11).Value = 6438.1665
$ws.Cells.Item(7, 13).Value = -6326.1665
$ws.Cells.Item(46, 8).Value = 3297.5
$ws.Cells.Item(46, 9).Value = 2198.75
$ws.Cells.Item(46, 10).Value = 5495
$ws.Cells.Item(46, 11).Value = 2198.75
$ws.Cells.Item(46, 12).Value = 5495
$ws.Cells.Item(46, 13).Value = -2010.75
$ws.Cells.Item(46, 14).Value = -5871
$ws.Cells.Item(122, 8).Value = 2792.3333
$ws.Cells.Item(122, 9).Value = 2579.8
$ws.Cells.Item(122, 11).Value = 7739.400000000001
$ws.Cells.Item(122, 13).Value = -5289.400000000001
$ws.Cells.Item(126, 8).Value = 6506
$ws.Cells.Item(126, 9).Value = 6438.1665
$ws.Cells.Item(126, 11).Value = 19314.4995
$ws.Cells.Item(126, 13).Value = -16844.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 4202.5386
$ws.Cells.Item(122, 9).Value = 3143.3
$ws.Cells.Item(122, 10).Value = 7733.3335
$ws.Cells.Item(122, 11).Value = 9429.900000000001
$ws.Cells.Item(122, 12).Value = 23200.0005
$ws.Cells.Item(122, 13).Value = -6979.900000000001
$ws.Cells.Item(122, 14).Value = -28100.0005
